$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1838
$ws1.Range("F6").Value = 475
$ws1.Range("F9").Value = 2440
$ws1.Range("F10").Value = 145
$ws1.Range("F11").Value = 79
$ws1.Range("F12").Value = 164
$ws1.Range("F13").Value = 1480
$ws1.Range("F14").Value = 517
$ws1.Range("F16").Value = 318
$ws1.Range("F20").Value = 201
$ws1.Range("F24").Value = 125
$ws1.Range("F26").Value = 1519
$ws1.Range("F29").Value = 352
$ws1.Range("F32").Value = 384

# Sheet "全部类型" (fourth sheet) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1838
$ws4.Range("F7").Value = 475
$ws4.Range("F10").Value = 2441
$ws4.Range("F11").Value = 145
$ws4.Range("F12").Value = 79
$ws4.Range("F13").Value = 164
$ws4.Range("F14").Value = 1480
$ws4.Range("F15").Value = 517
$ws4.Range("F17").Value = 318
$ws4.Range("F21").Value = 201
$ws4.Range("F25").Value = 125
$ws4.Range("F27").Value = 1519
$ws4.Range("F30").Value = 352
$ws4.Range("F33").Value = 384

$wb.Save()
